$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 29: LeetCode 69 - Sqrt(x)
$ws.Range("A29").Value = 69
$ws.Range("B29").Value = "Sqrt(x)"
$ws.Range("C29").Value = "#math #binary-search #重点 "
$ws.Range("D29").Value = "easy"
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 45838
$ws.Range("I29").Value = 45838

$ws.Range("A29:I29").HorizontalAlignment = -4108
$ws.Range("A29:I29").VerticalAlignment = -4108
$ws.Range("B29:C29").WrapText = $true

# Reuse the existing date style (s=3) instead of letting NumberFormat mint a new one
$ws.Range("H28:I28").Copy()
$ws.Range("H29:I29").PasteSpecial(-4122)
$ws.Range("H29").Value = 45838
$ws.Range("I29").Value = 45838

$ws.Rows.Item(29).RowHeight = 34

# update selection to mimic the author's later navigation
$ws.Range("E35").Select()
